# Auto-generated edit script applying value updates per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 30
$ws.Range("H33").Value = 594.5357
$ws.Range("I33").Value = 593.13635
$ws.Range("K33").Value = 593.13635
$ws.Range("M33").Value = -364.13635
$ws.Range("H132").Value = 2201.0334
$ws.Range("I132").Value = 1862.88
$ws.Range("K132").Value = 5588.64
$ws.Range("M132").Value = -3058.64
$ws.Range("H135").Value = 370942.4
$ws.Range("I135").Value = 417160.25
$ws.Range("J135").Value = 1199.6666
$ws.Range("K135").Value = 3754442.25
$ws.Range("L135").Value = 10796.9994
$ws.Range("M135").Value = -3751907.25
$ws.Range("N135").Value = -15866.9994
$ws.Range("H136").Value = 62500
$ws.Range("I136").Value = 50000
$ws.Range("J136").Value = 75000
$ws.Range("K136").Value = 50000
$ws.Range("L136").Value = 75000
$ws.Range("M136").Value = -44900
$ws.Range("N136").Value = -85200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 58527.137
$ws.Range("J74").Value = 4988.7
$ws.Range("L74").Value = 4988.7
$ws.Range("N74").Value = -6736.7
$ws.Range("H77").Value = 58527.137
$ws.Range("J77").Value = 4988.7
$ws.Range("L77").Value = 24943.5
$ws.Range("N77").Value = -33679.5
$ws.Range("H102").Value = 1035.72
$ws.Range("I102").Value = 1028.65
$ws.Range("J102").Value = 1064
$ws.Range("K102").Value = 1028.65
$ws.Range("L102").Value = 1064
$ws.Range("M102").Value = 593.3499999999999
$ws.Range("N102").Value = -4308
$ws.Range("H122").Value = 12490.523
$ws.Range("I122").Value = 15589.571
$ws.Range("J122").Value = 6292.4287
$ws.Range("K122").Value = 46768.713
$ws.Range("L122").Value = 18877.2861
$ws.Range("M122").Value = -44318.713
$ws.Range("N122").Value = -23777.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 69829.75
$ws.Range("J140").Value = 69829.75
$ws.Range("L140").Value = 69829.75
$ws.Range("N140").Value = -80189.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9509.82
$ws.Range("I31").Value = 4513.875
$ws.Range("J31").Value = 12985.261
$ws.Range("K31").Value = 4513.875
$ws.Range("L31").Value = 12985.261
$ws.Range("M31").Value = -4218.875
$ws.Range("N31").Value = -13575.261
$ws.Range("H34").Value = 9509.82
$ws.Range("I34").Value = 4513.875
$ws.Range("J34").Value = 12985.261
$ws.Range("K34").Value = 4513.875
$ws.Range("L34").Value = 12985.261
$ws.Range("M34").Value = -4311.875
$ws.Range("N34").Value = -13389.261

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3080.6956
$ws.Range("J34").Value = 5210.385
$ws.Range("L34").Value = 15631.155
$ws.Range("N34").Value = -15799.155
$ws.Range("H44").Value = 703.4
$ws.Range("I44").Value = 260.8
$ws.Range("K44").Value = 782.4000000000001
$ws.Range("M44").Value = -384.4000000000001
$ws.Range("H56").Value = 5432.5
$ws.Range("I56").Value = 5432.5
$ws.Range("K56").Value = 5432.5
$ws.Range("M56").Value = -4902.5
$ws.Range("H118").Value = 499
$ws.Range("I118").Value = 499
$ws.Range("K118").Value = 1497
$ws.Range("M118").Value = -254
$ws.Range("H132").Value = 13411.529
$ws.Range("I132").Value = 5777.6665
$ws.Range("J132").Value = 21999.625
$ws.Range("K132").Value = 51998.9985
$ws.Range("L132").Value = 197996.625
$ws.Range("M132").Value = -49468.9985
$ws.Range("N132").Value = -203056.625
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H134").Value = 70681.53
$ws.Range("I134").Value = 70681.53
$ws.Range("K134").Value = 212044.59
$ws.Range("M134").Value = -206974.59
$ws.Range("H136").Value = 714.2857
$ws.Range("I136").Value = 714.2857
$ws.Range("K136").Value = 2142.8571
$ws.Range("M136").Value = 2957.1429
$ws.Range("H137").Value = 125862.69
$ws.Range("I137").Value = 67586.92999999999
$ws.Range("K137").Value = 202760.79
$ws.Range("M137").Value = -197660.79
$ws.Range("H138").Value = 60470.723
$ws.Range("I138").Value = 75246.14
$ws.Range("K138").Value = 225738.42
$ws.Range("M138").Value = -220598.42
$ws.Range("H139").Value = 67044
$ws.Range("I139").Value = 503749.5
$ws.Range("K139").Value = 1511248.5
$ws.Range("M139").Value = -1506108.5
$ws.Range("H140").Value = 252115.88
$ws.Range("I140").Value = 334654.66
$ws.Range("K140").Value = 1003963.98
$ws.Range("M140").Value = -998783.98
$ws.Range("H141").Value = 11074.8
$ws.Range("I141").Value = 2687.25
$ws.Range("K141").Value = 8061.75
$ws.Range("M141").Value = -2881.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10169.75
$ws.Range("I70").Value = 9271.75
$ws.Range("K70").Value = 9271.75
$ws.Range("M70").Value = -9001.75
$ws.Range("H73").Value = 10169.75
$ws.Range("I73").Value = 9271.75
$ws.Range("K73").Value = 9271.75
$ws.Range("M73").Value = -8335.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5330.1
$ws.Range("I7").Value = 2963.2727
$ws.Range("K7").Value = 2963.2727
$ws.Range("M7").Value = -2851.2727
$ws.Range("H16").Value = 7799.6
$ws.Range("J16").Value = 7999.6665
$ws.Range("L16").Value = 7999.6665
$ws.Range("N16").Value = -8339.666499999999
$ws.Range("H22").Value = 3321.5386
$ws.Range("I22").Value = 2033.3334
$ws.Range("J22").Value = 3708
$ws.Range("K22").Value = 2033.3334
$ws.Range("L22").Value = 3708
$ws.Range("M22").Value = -1738.3334
$ws.Range("N22").Value = -4298
$ws.Range("H27").Value = 3321.5386
$ws.Range("I27").Value = 2033.3334
$ws.Range("J27").Value = 3708
$ws.Range("K27").Value = 2033.3334
$ws.Range("L27").Value = 3708
$ws.Range("M27").Value = -1926.3334
$ws.Range("N27").Value = -3922
$ws.Range("H46").Value = 2009.1333
$ws.Range("I46").Value = 458.6
$ws.Range("K46").Value = 458.6
$ws.Range("M46").Value = -270.6
$ws.Range("I76").Value = 29999
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 29999
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -29661
$ws.Range("N76").ClearContents()
$ws.Range("I79").Value = 29999
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 29999
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -28829
$ws.Range("N79").ClearContents()
$ws.Range("H126").Value = 5330.1
$ws.Range("I126").Value = 2963.2727
$ws.Range("K126").Value = 8889.8181
$ws.Range("M126").Value = -6419.8181
$ws.Range("H136").Value = 9811.581
$ws.Range("J136").Value = 12260.435
$ws.Range("L136").Value = 36781.305
$ws.Range("N136").Value = -41881.305
